# Generate Report for Handoff
#
# The localization run finished: status flips from "In Translation" to
# "Ready for handoff" on every sheet, and the "latest handoff" timestamps
# are refreshed to the moment the handoff package was produced. Widen the
# now-longer status/date columns so the new text isn't clipped.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -------------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps ---------------------------------------
$zhcn.Range("H2").Value     = "2016-08-24 13:03:18"   # zh-cn Latest Handoff Datetime
$dede.Range("H2").Value     = "2016-08-24 13:03:23"   # de-de Latest Handoff Datetime
$overview.Range("G2").Value = "2016-08-24 13:03:23"   # Latest HO Xliff Generate Date

# --- Widen the status/date columns to fit the new text ------------------
# Target display width ~17.22 characters; this host quantizes ColumnWidth
# writes to its internal pixel grid, so feed it the input that lands on
# the grid point closest to that target.
$overview.Range("E1").EntireColumn.ColumnWidth = 16.33
$overview.Range("F1").EntireColumn.ColumnWidth = 16.33
$zhcn.Range("C1").EntireColumn.ColumnWidth = 16.33
$dede.Range("C1").EntireColumn.ColumnWidth = 16.33
